# Fix accidental re-use of CKR: the Cow IDs (and the GaitScore readings
# recorded against them) in rows 2-51 were mixed up with IDs already
# used elsewhere; replace them with the correct, non-duplicated IDs/scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values are 5-digit Cow IDs that can have leading zeros, so they
# must stay text. Force Text format before writing the new IDs, then drop
# back to the default style so no stray number format is left behind.
$colA = $ws.Range("A2:A51")
$colA.NumberFormat = "@"

$ws.Range("A2").Value = "76490"
$ws.Range("A3").Value = "27588"
$ws.Range("A4").Value = "04438"
$ws.Range("A5").Value = "43840"
$ws.Range("A6").Value = "84609"
$ws.Range("A7").Value = "36363"
$ws.Range("A8").Value = "25391"
$ws.Range("A9").Value = "95016"
$ws.Range("A10").Value = "91689"
$ws.Range("A11").Value = "12561"
$ws.Range("A12").Value = "93225"
$ws.Range("A13").Value = "82422"
$ws.Range("A14").Value = "49032"
$ws.Range("A15").Value = "21352"
$ws.Range("A16").Value = "22313"
$ws.Range("A17").Value = "73246"
$ws.Range("A18").Value = "91991"
$ws.Range("A19").Value = "97764"
$ws.Range("A20").Value = "50859"
$ws.Range("A21").Value = "15092"
$ws.Range("A22").Value = "93744"
$ws.Range("A23").Value = "48412"
$ws.Range("A24").Value = "70366"
$ws.Range("A25").Value = "36418"
$ws.Range("A26").Value = "03934"
$ws.Range("A27").Value = "39143"
$ws.Range("A28").Value = "79916"
$ws.Range("A29").Value = "22707"
$ws.Range("A30").Value = "72769"
$ws.Range("A31").Value = "96038"
$ws.Range("A32").Value = "23910"
$ws.Range("A33").Value = "98581"
$ws.Range("A34").Value = "51246"
$ws.Range("A35").Value = "39597"
$ws.Range("A36").Value = "20428"
$ws.Range("A37").Value = "22117"
$ws.Range("A38").Value = "57701"
$ws.Range("A39").Value = "21573"
$ws.Range("A40").Value = "65391"
$ws.Range("A41").Value = "72985"
$ws.Range("A42").Value = "90541"
$ws.Range("A43").Value = "38417"
$ws.Range("A44").Value = "25884"
$ws.Range("A45").Value = "10974"
$ws.Range("A46").Value = "14473"
$ws.Range("A47").Value = "79964"
$ws.Range("A48").Value = "90530"
$ws.Range("A49").Value = "26743"
$ws.Range("A50").Value = "19258"
$ws.Range("A51").Value = "57153"

$colA.Style = "Normal"

# Row 43 previously held a free-text Note (no GaitScore); it now gets a
# normal GaitScore of 0 and loses the note.
$ws.Range("C43").ClearContents()
$ws.Range("B43").Value = 0

# Row 46 previously held a normal GaitScore of 1; it now instead gets a
# free-text Note explaining why no GaitScore was recorded.
$ws.Range("B46").ClearContents()
$ws.Range("C46").Value = "Not being milked due to clinical mastitis"

$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 1
$ws.Range("B24").Value = 2
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 1
$ws.Range("B31").Value = 0
$ws.Range("B32").Value = 0
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = 0
$ws.Range("B35").Value = 0
$ws.Range("B36").Value = 0
$ws.Range("B37").Value = 0
$ws.Range("B38").Value = 0
$ws.Range("B39").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("B42").Value = 0
$ws.Range("B44").Value = 0
$ws.Range("B45").Value = 0
$ws.Range("B47").Value = 0
$ws.Range("B48").Value = 0
$ws.Range("B49").Value = 0
$ws.Range("B50").Value = 1
$ws.Range("B51").Value = 0
